$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header columns for the log template dropdown (Top3, Top4, Top5)
$ws.Range("E1").Value = "Top3"
$ws.Range("F1").Value = "Top4"
$ws.Range("G1").Value = "Top5"

# Update existing depth values and add new columns' data for row 2
$ws.Range("C2").Value = 400
$ws.Range("D2").Value = 600
$ws.Range("E2").Value = 800
$ws.Range("F2").Value = 1000
$ws.Range("G2").Value = 1250

# Adjust column widths to match the new, narrower UI layout
$ws.Range("A1").ColumnWidth = 11.33
$ws.Range("B1").ColumnWidth = 4.17
$ws.Range("C1:G1").ColumnWidth = 4.67

# Reflect the new active cell/selection after the edits
$ws.Range("G2").Select()
